$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for all data rows (2-472)
# from 45182 to 45184 (Excel serial date numbers).
$ws.Range("C2:C472").Value = 45184
